$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the literal string into the cell without Excel's
    # numeric auto-detection coercing it to a Double (e.g. "1.001" -> 1.001).
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "29.786.69"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "1.888.23"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  -5.54%  "
Set-TextCell "D6" "244.45"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -4.25%  "
Set-TextCell "D9" "25.19"
$ws.Range("E9").Value = "  -7.82%  "
Set-TextCell "D10" "0.07215"
$ws.Range("E10").Value = "  -1.00%  "
Set-TextCell "D11" "0.08082"
$ws.Range("E11").Value = "  -0.37%  "
Set-TextCell "D12" "0.7643"
$ws.Range("E12").Value = "  -4.05%  "
Set-TextCell "D13" "5.480"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "1.891.04"
$ws.Range("E14").Value = "  -2.76%  "
Set-TextCell "D15" "92.24"
$ws.Range("E15").Value = "  -2.24%  "
Set-TextCell "D16" "6.177"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "29.791.53"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("E18").Value = "  -2.75%  "
Set-TextCell "D19" "242.33"
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D21" "1.001"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.154.78"
$ws.Range("E22").Value = "  -1.11%  "
Set-TextCell "D23" "8.147"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("E24").Value = "  +0.02%  "
Set-TextCell "D25" "0.1579"
$ws.Range("E25").Value = "  -6.12%  "
$ws.Range("E26").Value = "  -1.00%  "
Set-TextCell "D27" "162.50"
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("E28").Value = "  -1.83%  "
Set-TextCell "D29" "2.034"
$ws.Range("E29").Value = "  -5.65%  "
Set-TextCell "D30" "1.436"
$ws.Range("E30").Value = "  +4.73%  "
Set-TextCell "D31" "1.548"
$ws.Range("E31").Value = "  -0.27%  "
Set-TextCell "D32" "4.453"
$ws.Range("E32").Value = "  +2.28%  "
Set-TextCell "D33" "4.074"
$ws.Range("E33").Value = "  -1.82%  "
Set-TextCell "D34" "0.05493"
$ws.Range("E34").Value = "  -3.79%  "
Set-TextCell "D35" "1.256"
$ws.Range("E35").Value = "  -3.65%  "
Set-TextCell "D36" "0.7486"
$ws.Range("E36").Value = "  +0.01%  "
Set-TextCell "D37" "0.9994"
$ws.Range("E37").Value = "  -0.01%  "
Set-TextCell "D38" "2.633"
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("E39").Value = "  -2.10%  "
Set-TextCell "D40" "2.782"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("D41").Value = "1.151.05"
$ws.Range("E41").Value = "  +10.76%  "
Set-TextCell "D42" "73.60"
$ws.Range("E42").Value = "  -1.86%  "
Set-TextCell "D43" "0.4416"
$ws.Range("E43").Value = "  -2.11%  "
Set-TextCell "D44" "5.892"
$ws.Range("E44").Value = "  -1.71%  "
Set-TextCell "D45" "0.8481"
$ws.Range("E45").Value = "  -0.85%  "
Set-TextCell "D46" "1.001"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D47" "102.58"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D48" "1.879"
$ws.Range("E48").Value = "  -2.63%  "
Set-TextCell "D49" "9.989"
$ws.Range("E49").Value = "  +0.32%  "
Set-TextCell "D50" "7.436"
$ws.Range("E50").Value = "  -2.99%  "
Set-TextCell "D51" "3.013"
$ws.Range("E51").Value = "  -2.73%  "
